$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that can differ between the two matches sharing a kickoff slot:
# B (id) and F..AC (HomeTeam through PL_AhUnder). A/C/D/E (row idx/Div/Div
# Original Name/Date) are identical for swapped rows and are left untouched.
$cols = @(2) + @(6..29)

# 2023/24 Finland Division 1 results were re-matched to the correct fixtures;
# snapshot the current (pre-fix) row contents first so the swaps below do not
# clobber each other, then write the corrected values back out.
$rows = @(12, 16, 24, 25, 26, 28, 29, 59, 60, 99, 100, 101, 102, 110, 111, 113, 114, 118, 119, 120)
$snapshot = @{}
foreach ($r in $rows) {
    $rowData = @{}
    foreach ($c in $cols) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowData
}

# target row -> source row whose data it should now hold
$mapping = @{
    12 = 16
    16 = 12
    24 = 26
    25 = 24
    26 = 25
    28 = 29
    29 = 28
    59 = 60
    60 = 59
    99 = 100
    100 = 101
    101 = 102
    102 = 99
    110 = 111
    111 = 110
    113 = 114
    114 = 113
    118 = 119
    119 = 120
    120 = 118
}

foreach ($target in $mapping.Keys) {
    $source = $mapping[$target]
    $sourceData = $snapshot[$source]
    foreach ($c in $cols) {
        $ws.Cells.Item($target, $c).Value = $sourceData[$c]
    }
}
